$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 9; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 14; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 33; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 48; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 54; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 56; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 63; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 72; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 75; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 81; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 83; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 85; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 95; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 99; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 106; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 108; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 111; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 122; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 124; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 126; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 132; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 147; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 162; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 167; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 190; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 191; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 193; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 197; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 198; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 201; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 208; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 209; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 219; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 232; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 233; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 238; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 240; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 243; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 244; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 246; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 248; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 249; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 253; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 254; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 261; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 277; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 284; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()
